$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PCB_C")

# Fix BOM quantities: J3 (row 10) 1 -> 2, J4 (row 11) 1 -> 4
$ws.Range("A10").Value = 2
$ws.Range("A11").Value = 4

# Update the last active selection to match the saved file
$ws.Range("C23").Select()
